# Generated PowerShell COM-interop script
# Implements the ADANIPORTS.NS.xlsx sheet1 edit described by the diff:
#  1. Zero out column R for a set of historical rows (dividend/split backup flag reset)
#  2. Row 61 also has column Q reset to 0
#  3. Row 864: column O (isPivot) changes from 0 to 3
#  4. Rows 866-867: column R becomes an explicit numeric 0 (was blank)
#  5. Nine new weekly rows (868-876) are appended with OHLCV + derived columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: reset column R (and Q61) back to 0 for rows that previously carried a
#            nonzero 'backup'/'detect_structure' marker ---
$rowsToZeroR = @(61, 66, 77, 88, 92, 108, 119, 130, 137, 157, 164, 170, 176, 178, 183, 187, 191, 201, 212, 216, 222, 229, 241, 254, 255, 272, 277, 279, 285, 294, 297, 304, 316, 323, 330, 343, 360, 367, 375, 381, 384, 436, 443, 464, 473, 477, 514, 521, 527, 533, 538, 546, 553, 560, 566, 575, 620, 635, 651, 665, 669, 674, 677, 705, 714, 715, 721, 733, 739, 745, 760, 781, 800, 805, 816, 832, 838, 858, 861)
foreach ($r in $rowsToZeroR) {
    $ws.Cells.Item($r, 18).Value = 0
}

# Row 61 additionally had column Q (detect_structure) reset to 0
$ws.Cells.Item(61, 17).Value = 0

# --- 3: row 864, column O (isPivot) 0 -> 3 ---
$ws.Cells.Item(864, 15).Value = 3

# --- 4: rows 866 & 867, column R blank -> numeric 0 ---
$ws.Cells.Item(866, 18).Value = 0
$ws.Cells.Item(867, 18).Value = 0

# --- 5: append new weekly rows 868-876 ---
# Row 868
$ws.Cells.Item(868, 1).Value = 45474
$ws.Cells.Item(868, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(868, 2).Value = 1483
$ws.Cells.Item(868, 3).Value = 1520
$ws.Cells.Item(868, 4).Value = 1455.050048828125
$ws.Cells.Item(868, 5).Value = 1500.449951171875
$ws.Cells.Item(868, 6).Value = 1500.449951171875
$ws.Cells.Item(868, 7).Value = 18269328
$ws.Cells.Item(868, 8).Value = 2024
$ws.Cells.Item(868, 9).Value = 7
$ws.Cells.Item(868, 10).Value = 1
$ws.Cells.Item(868, 11).Value = 0
$ws.Cells.Item(868, 12).Value = 0
$ws.Cells.Item(868, 13).Value = 0
$ws.Cells.Item(868, 14).Value = 27
$ws.Cells.Item(868, 15).Value = 0
$ws.Cells.Item(868, 16).Value = 0
$ws.Cells.Item(868, 17).Value = 0

# Row 869
$ws.Cells.Item(869, 1).Value = 45481
$ws.Cells.Item(869, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(869, 2).Value = 1496
$ws.Cells.Item(869, 3).Value = 1503.949951171875
$ws.Cells.Item(869, 4).Value = 1460
$ws.Cells.Item(869, 5).Value = 1486.699951171875
$ws.Cells.Item(869, 6).Value = 1486.699951171875
$ws.Cells.Item(869, 7).Value = 10300965
$ws.Cells.Item(869, 8).Value = 2024
$ws.Cells.Item(869, 9).Value = 7
$ws.Cells.Item(869, 10).Value = 8
$ws.Cells.Item(869, 11).Value = 0
$ws.Cells.Item(869, 12).Value = 0
$ws.Cells.Item(869, 13).Value = 0
$ws.Cells.Item(869, 14).Value = 28
$ws.Cells.Item(869, 15).Value = 0
$ws.Cells.Item(869, 16).Value = 0
$ws.Cells.Item(869, 17).Value = 0

# Row 870
$ws.Cells.Item(870, 1).Value = 45488
$ws.Cells.Item(870, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(870, 2).Value = 1494
$ws.Cells.Item(870, 3).Value = 1510.599975585938
$ws.Cells.Item(870, 4).Value = 1465.5
$ws.Cells.Item(870, 5).Value = 1469.300048828125
$ws.Cells.Item(870, 6).Value = 1469.300048828125
$ws.Cells.Item(870, 7).Value = 9915071
$ws.Cells.Item(870, 8).Value = 2024
$ws.Cells.Item(870, 9).Value = 7
$ws.Cells.Item(870, 10).Value = 15
$ws.Cells.Item(870, 11).Value = 0
$ws.Cells.Item(870, 12).Value = 0
$ws.Cells.Item(870, 13).Value = 0
$ws.Cells.Item(870, 14).Value = 29
$ws.Cells.Item(870, 15).Value = 0
$ws.Cells.Item(870, 16).Value = 0
$ws.Cells.Item(870, 17).Value = 0

# Row 871
$ws.Cells.Item(871, 1).Value = 45495
$ws.Cells.Item(871, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(871, 2).Value = 1470
$ws.Cells.Item(871, 3).Value = 1547.900024414062
$ws.Cells.Item(871, 4).Value = 1419
$ws.Cells.Item(871, 5).Value = 1542.75
$ws.Cells.Item(871, 6).Value = 1542.75
$ws.Cells.Item(871, 7).Value = 16622452
$ws.Cells.Item(871, 8).Value = 2024
$ws.Cells.Item(871, 9).Value = 7
$ws.Cells.Item(871, 10).Value = 22
$ws.Cells.Item(871, 11).Value = 0
$ws.Cells.Item(871, 12).Value = 0
$ws.Cells.Item(871, 13).Value = 0
$ws.Cells.Item(871, 14).Value = 30
$ws.Cells.Item(871, 15).Value = 2
$ws.Cells.Item(871, 16).Value = 0
$ws.Cells.Item(871, 17).Value = 0

# Row 872
$ws.Cells.Item(872, 1).Value = 45502
$ws.Cells.Item(872, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(872, 2).Value = 1547.900024414062
$ws.Cells.Item(872, 3).Value = 1604.949951171875
$ws.Cells.Item(872, 4).Value = 1536
$ws.Cells.Item(872, 5).Value = 1588
$ws.Cells.Item(872, 6).Value = 1588
$ws.Cells.Item(872, 7).Value = 18825426
$ws.Cells.Item(872, 8).Value = 2024
$ws.Cells.Item(872, 9).Value = 7
$ws.Cells.Item(872, 10).Value = 29
$ws.Cells.Item(872, 11).Value = 0
$ws.Cells.Item(872, 12).Value = 0
$ws.Cells.Item(872, 13).Value = 0
$ws.Cells.Item(872, 14).Value = 31
$ws.Cells.Item(872, 15).Value = 0
$ws.Cells.Item(872, 16).Value = 0
$ws.Cells.Item(872, 17).Value = 0

# Row 873
$ws.Cells.Item(873, 1).Value = 45509
$ws.Cells.Item(873, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(873, 2).Value = 1540
$ws.Cells.Item(873, 3).Value = 1563.449951171875
$ws.Cells.Item(873, 4).Value = 1471.550048828125
$ws.Cells.Item(873, 5).Value = 1533.800048828125
$ws.Cells.Item(873, 6).Value = 1533.800048828125
$ws.Cells.Item(873, 7).Value = 16536798
$ws.Cells.Item(873, 8).Value = 2024
$ws.Cells.Item(873, 9).Value = 8
$ws.Cells.Item(873, 10).Value = 5
$ws.Cells.Item(873, 11).Value = 0
$ws.Cells.Item(873, 12).Value = 0
$ws.Cells.Item(873, 13).Value = 0
$ws.Cells.Item(873, 14).Value = 32
$ws.Cells.Item(873, 15).Value = 0
$ws.Cells.Item(873, 16).Value = 0
$ws.Cells.Item(873, 17).Value = 0

# Row 874
$ws.Cells.Item(874, 1).Value = 45516
$ws.Cells.Item(874, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(874, 2).Value = 1500
$ws.Cells.Item(874, 3).Value = 1524.25
$ws.Cells.Item(874, 4).Value = 1452
$ws.Cells.Item(874, 5).Value = 1493.449951171875
$ws.Cells.Item(874, 6).Value = 1493.449951171875
$ws.Cells.Item(874, 7).Value = 16349330
$ws.Cells.Item(874, 8).Value = 2024
$ws.Cells.Item(874, 9).Value = 8
$ws.Cells.Item(874, 10).Value = 12
$ws.Cells.Item(874, 11).Value = 0
$ws.Cells.Item(874, 12).Value = 0
$ws.Cells.Item(874, 13).Value = 0
$ws.Cells.Item(874, 14).Value = 33
$ws.Cells.Item(874, 15).Value = 0
$ws.Cells.Item(874, 16).Value = 0
$ws.Cells.Item(874, 17).Value = 0

# Row 875
$ws.Cells.Item(875, 1).Value = 45523
$ws.Cells.Item(875, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(875, 2).Value = 1500.949951171875
$ws.Cells.Item(875, 3).Value = 1517
$ws.Cells.Item(875, 4).Value = 1486.099975585938
$ws.Cells.Item(875, 5).Value = 1491.300048828125
$ws.Cells.Item(875, 6).Value = 1491.300048828125
$ws.Cells.Item(875, 7).Value = 6508605
$ws.Cells.Item(875, 8).Value = 2024
$ws.Cells.Item(875, 9).Value = 8
$ws.Cells.Item(875, 10).Value = 19
$ws.Cells.Item(875, 11).Value = 0
$ws.Cells.Item(875, 12).Value = 0
$ws.Cells.Item(875, 13).Value = 0
$ws.Cells.Item(875, 14).Value = 34
$ws.Cells.Item(875, 15).Value = 0
$ws.Cells.Item(875, 16).Value = 0
$ws.Cells.Item(875, 17).Value = 0

# Row 876
$ws.Cells.Item(876, 1).Value = 45530
$ws.Cells.Item(876, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(876, 2).Value = 1495.5
$ws.Cells.Item(876, 3).Value = 1497
$ws.Cells.Item(876, 4).Value = 1452.599975585938
$ws.Cells.Item(876, 5).Value = 1481.900024414062
$ws.Cells.Item(876, 6).Value = 1481.900024414062
$ws.Cells.Item(876, 7).Value = 8578184
$ws.Cells.Item(876, 8).Value = 2024
$ws.Cells.Item(876, 9).Value = 8
$ws.Cells.Item(876, 10).Value = 26
$ws.Cells.Item(876, 11).Value = 0
$ws.Cells.Item(876, 12).Value = 0
$ws.Cells.Item(876, 13).Value = 0
$ws.Cells.Item(876, 14).Value = 35
$ws.Cells.Item(876, 15).Value = 0
$ws.Cells.Item(876, 16).Value = 0
$ws.Cells.Item(876, 17).Value = 0

